$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The "LimType" bound-type label used for the NCAP_BND table was renamed to
# CAP_BND; update every row of the block (C7:C26) that carried the old label.
$ws.Range("C7:C26").Value = "CAP_BND"

# Scroll the view down so row 13 is the top visible row, then leave the
# selection on the block of cells that was just edited.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C7:C26").Select()
